# TradingModel - 2021/11/16 data update
# Append three new trading-history rows (2021/11/16) to the bottom of the
# existing log on the single worksheet, then move the view/selection down
# to show the newly entered rows (mirrors Excel scrolling to F31 after
# typing the last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateSerial   = 44516          # 2021/11/16
$dateFormat   = "m""月""d""日"""

# Date, Stock_Id, Action, PositionSize, Price
$newRows = @(
    @(29, 8289, "long", 160, 36.85),
    @(30, 2314, "long",  65, 93),
    @(31, 3221, "long",  61, 42.6)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $dateSerial
    $dateCell.NumberFormat = $dateFormat

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Scroll the view down to the new rows and leave the selection where Excel
# would land after entering the last value (one column past the data).
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F31").Select()
